# Update profit files after running on 2025-10-10
# Append a new row (row 54) to the profit log: date 10/10/2025, profit 14519.98

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date cell to be stored as plain text (matching the existing rows,
# which store dates as literal strings like "10/09/2025" rather than real
# date values) instead of letting Excel auto-convert the string to a date.
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = "10/10/2025"
# Reset the cell style back to the default "Normal" style so no explicit
# number-format/style index is left behind on the cell (matching the rest
# of the sheet, which uses the default style).
$ws.Range("A54").Style = "Normal"

$ws.Range("B54").Value = 14519.98
